$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Part 1: rewrite "Figure out what the pattern would be to predict  where she
# would end up at any multiple of 10. " -> split across 3 runs, drop the
# proofErr gramStart/gramEnd markers, and collapse the double space between
# "predict" and "where" into a single space.
# ---------------------------------------------------------------------------

$oldSpan = "Figure out what the pattern would be to " + "predict  where" + " she would end up at any multiple of 10. "
$piece2 = "Figure out what t"
$piece3 = "he pattern would be to predict "
$piece4 = "where she would end up at any multiple of 10. "

$full = $d.Content.Text
$spanStart = $full.IndexOf($oldSpan)
if ($spanStart -lt 0) { throw "could not locate target span" }

$r = $d.Range($spanStart, $spanStart + $oldSpan.Length)
$r.Text = $piece2 + $piece3 + $piece4

# Build clean (unformatted) donor copies of piece2/piece3/piece4 in a scratch
# paragraph appended at the end of the document; Range.FormattedText lets us
# drop replacement text into the real location as an independent run without
# leaving a stray (but empty) rPr behind the way Font property churn would.
$scratchPara = $d.Paragraphs.Add()
$scratchPara.Range.Text = $piece2 + "|" + $piece3 + "|" + $piece4
$scratchStart = $scratchPara.Range.Start

$donor2 = $d.Range($scratchStart, $scratchStart + $piece2.Length)
$donor3Start = $scratchStart + $piece2.Length + 1
$donor3 = $d.Range($donor3Start, $donor3Start + $piece3.Length)
$donor4Start = $donor3Start + $piece3.Length + 1
$donor4 = $d.Range($donor4Start, $donor4Start + $piece4.Length)

$ft2 = $donor2.FormattedText
$ft3 = $donor3.FormattedText
$ft4 = $donor4.FormattedText

# Stamp from right to left so earlier offsets stay valid.
$p4Start = $spanStart + $piece2.Length + $piece3.Length
$p4End = $p4Start + $piece4.Length
$target4 = $d.Range($p4Start, $p4End)
$target4.FormattedText = $ft4

$p3Start = $spanStart + $piece2.Length
$p3End = $p3Start + $piece3.Length
$target3 = $d.Range($p3Start, $p3End)
$target3.FormattedText = $ft3

$p2Start = $spanStart
$p2End = $p2Start + $piece2.Length
$target2 = $d.Range($p2Start, $p2End)
$target2.FormattedText = $ft2

# Drop the scratch paragraph now that all donors have been consumed.
$scratchPara2 = $d.Paragraphs.Item($d.Paragraphs.Count)
$scratchPara2.Range.Delete()

# ---------------------------------------------------------------------------
# Part 2: add the new sentence as its own run in the (previously empty) list
# item paragraph, right before the _GoBack bookmark.
# ---------------------------------------------------------------------------

$apos = [char]0x2019
$newSentence = "The solution for this problem would be that the little girl would be ending on 10" + $apos + "s on her first finger and her ring finger in an alternating manner. After 10 it would alternate every 2 times, so 10 would be on the first finger, 20 and 30 would be on the ring finger and 40 and 50 would be on the first finger and so on up the number tree. "

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$insertPos = $lastPara.Range.Start
$insertRange = $d.Range($insertPos, $insertPos)
$insertRange.InsertBefore($newSentence)
